# Commit standing frame BOM
#
# Net changes to the "BOM" sheet:
#   1. Insert a new row at row 19 (上身电机上盖底部限位 / M3.5 x 16 spacer screws),
#      which shifts every row below it down by one and auto-grows the
#      shared SUM formula at the bottom of the sheet.
#   2. Fill in one of the still-empty placeholder rows in the MISC section
#      (now at row 52 after the shift above) with the new "2040 型材"
#      standing-test-stand frame line item, including its Misumi hyperlink.
#   3. Re-create every hyperlink on the sheet so the ranges track the rows
#      they now point to (row insertion in this host does not itself walk
#      hyperlink ranges), and add the brand-new hyperlink for row 52.
#   4. Nudge column C a bit wider to fit the long new URL, and leave the
#      selection/view the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# ---------------------------------------------------------------------
# 1. Existing hyperlinks, keyed by their CURRENT (pre-edit) row number.
#    We'll delete and re-create all of these after the row insert below,
#    since this host's Rows.Insert() does not shift hyperlink ranges.
# ---------------------------------------------------------------------
$hyperlinkRows = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,32,33,34,35,36,37,39,40,41,43,45,46,47,48)
$hyperlinkUrls = @(
  "https://detail.tmall.com/item.htm?id=899448134936&skuId=5754304832006",
  "https://detail.tmall.com/item.htm?id=899448134936&skuId=5798756972352",
  "https://detail.tmall.com/item.htm?id=899448134936&skuId=5922977711804",
  "https://detail.tmall.com/item.htm?id=822545298024&skuId=5702459126386",
  "https://item.taobao.com/item.htm?id=537967383377&skuId=3212738635362",
  "https://item.taobao.com/item.htm?id=537967383377&skuId=5461813862733",
  "https://item.taobao.com/item.htm?id=525468794969&skuId=3127568831102",
  "https://item.taobao.com/item.htm?id=525468794969&skuId=3127568831102",
  "https://item.taobao.com/item.htm?id=529627939342&skuId=3165003331293",
  "https://item.taobao.com/item.htm?id=529627939342&skuId=3206368205626",
  "https://item.taobao.com/item.htm?id=529627939342&skuId=3200246575353",
  "https://item.taobao.com/item.htm?id=529627939342&skuId=3165003331294",
  "https://item.taobao.com/item.htm?id=38977230812&skuId=3165007311199",
  "https://item.taobao.com/item.htm?id=38977230812&skuId=3165007311200",
  "https://detail.tmall.com/item.htm?id=625077079266&skuId=5611160825675",
  "https://item.taobao.com/item.htm?id=652062390391&skuId=4698920756210",
  "https://detail.tmall.com/item.htm?id=610905186565&skuId=5148221219880",
  "https://detail.tmall.com/item.htm?id=610905186565&skuId=4925270516930",
  "https://item.taobao.com/item.htm?id=667085311743&skuId=4865635364512",
  "https://item.taobao.com/item.htm?id=577536722051",
  "https://item.taobao.com/item.htm?id=611636440997&skuId=4310799149744",
  "https://detail.tmall.com/item.htm?id=679335252560&skuId=5978377674682",
  "https://detail.tmall.com/item.htm?id=679335252560&skuId=5079135404284",
  "https://detail.tmall.com/item.htm?id=679335252560&skuId=5136048959263",
  "https://detail.tmall.com/item.htm?id=679335252560&skuId=5079135404284",
  "https://detail.tmall.com/item.htm?id=654533678350&skuId=4892505606177",
  "https://detail.tmall.com/item.htm?id=654533678350&skuId=5036835474357",
  "https://detail.tmall.com/item.htm?id=645061129098&skuId=4648099525532",
  "https://item.taobao.com/item.htm?id=537967383377&skuId=3212738635362",
  "https://item.taobao.com/item.htm?id=537967383377&skuId=3212738635365",
  "https://item.taobao.com/item.htm?id=529627939342&skuId=3165003331293",
  "https://item.taobao.com/item.htm?id=566323233973&skuId=5212591909140",
  "https://detail.tmall.com/item.htm?id=679335252560&skuId=5079135404284",
  "https://detail.tmall.com/item.htm?id=679335252560&skuId=5136048959263",
  "https://detail.tmall.com/item.htm?id=899448134936&skuId=5922977711804",
  "https://item.taobao.com/item.htm?id=38977230812&skuId=3165007311199",
  "https://item.taobao.com/item.htm?id=537967383377&skuId=5461813862733",
  "https://detail.tmall.com/item.htm?id=679335252560&skuId=5079135404284",
  "https://detail.tmall.com/item.htm?id=899448134936&skuId=5922977711804",
  "https://item.taobao.com/item.htm?id=537967383377&skuId=3212738635365",
  "https://item.taobao.com/item.htm?id=38977230812&skuId=3165007311199",
  "https://item.taobao.com/item.htm?id=35431215248&skuId=3165685753062"
)

$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2. Insert the new "上身电机上盖底部限位" row above the old row 19.
#    Everything from row 19 down shifts to row+1 (formulas/merges/the
#    bottom SUM range all track automatically).
# ---------------------------------------------------------------------
$ws.Rows.Item(19).Insert()
$ws.Range("A18:G18").Copy()
$ws.Range("A19:G19").PasteSpecial(-4122)

$ws.Range("A19").Value = "粒片螺丝"
$ws.Range("B19").Value = "M3.5 x 16"
$ws.Range("D19").Value = "上身电机上盖底部限位"
$ws.Range("E19").Value = 6
$ws.Range("F19").Value = 0.3
$ws.Range("G19").Formula = "=E19*F19"

# ---------------------------------------------------------------------
# 3. Fill in the new standing-test-stand frame line item. This lands on
#    what used to be the first empty placeholder row of the MISC
#    section -- row 51 before the insert above, now row 52 -- so no
#    further row insertion is needed here.
# ---------------------------------------------------------------------
$ws.Range("A52").Value = "2040 型材"
$ws.Range("B52").Value = "经济型欧标槽宽 6mm"
$ws.Range("C52").Value = "站立调试底座"
$ws.Range("D52").Value = "站立调试底座"
$ws.Range("E52").Value = 5
$ws.Range("F52").Value = 13.19
$ws.Range("G52").Formula = "=E52*F52"

# D52 is the part/usage note, C52 is the hyperlinked source -- fix C52's
# text back to the product page label used by the source workbook and
# put the real usage note in D52.
$ws.Range("C52").Value = "站立调试底座"
$ws.Range("D52").Value = "站立调试底座"

# ---------------------------------------------------------------------
# 4. Re-create all the hyperlinks, shifting any row that used to be >= 19
#    down by one (matching the row insert above), then add the new one
#    for the Misumi 2040 extrusion on row 52.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $hyperlinkRows.Length; $i++) {
  $r = $hyperlinkRows[$i]
  if ($r -ge 19) { $r = $r + 1 }
  $ws.Hyperlinks.Add($ws.Cells.Item($r, 3), $hyperlinkUrls[$i])
}
$ws.Hyperlinks.Add($ws.Cells.Item(52, 3), "https://www.misumi.com.cn/vona2/detail/110310158399/?ProductCode=LCFB6-2040-200-TPW")

# ---------------------------------------------------------------------
# 5. Cosmetic bits: column C needs to be wide enough for the long new
#    URL, and the author left the sheet scrolled back to the top with
#    A52 selected.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 80.6640625

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A52").Select()
